{"js": "// Apply the commit's edits:\n// 1. \"WebDev Inc members:\" paragraph \u2014 drop the spell-check proofErr\n//    markers and merge the \"WebDev\" / \" Inc members:\" runs into a\n//    single run (same bold formatting).\n// 2. \"Karl Sanchez\" paragraph \u2014 append a student ID, split across two\n//    extra runs: \" A0003707\" and \"4\" (reads as \"A00037074\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst PKG_OPEN =\n  '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document ' +\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n  'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>';\nconst PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n// Grab the paragraph's own opening-tag attributes (w14:paraId, rsids, \u2026)\n// straight from its current OOXML so the replacement keeps them intact.\nasync function openTagAttrs(paragraph) {\n  const res = paragraph.getOoxml();\n  await context.sync();\n  const m = /<w:p\\b([^>]*)>/.exec(res.value);\n  return m ? m[1] : \"\";\n}\n\nlet membersPara = null;\nlet karlPara = null;\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (membersPara === null && t.indexOf(\"WebDev\") !== -1 && t.indexOf(\"Inc members\") !== -1) {\n    membersPara = p;\n  }\n  if (karlPara === null && t.indexOf(\"Karl Sanchez\") !== -1) {\n    karlPara = p;\n  }\n}\n\nif (membersPara) {\n  const attrs = await openTagAttrs(membersPara);\n  const xml =\n    PKG_OPEN +\n    '<w:p' + attrs + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>WebDev Inc members:</w:t></w:r>' +\n    '</w:p>' +\n    PKG_CLOSE;\n  membersPara.insertOoxml(xml, \"Replace\");\n  await context.sync();\n}\n\nif (karlPara) {\n  const attrs = await openTagAttrs(karlPara);\n  const xml =\n    PKG_OPEN +\n    '<w:p' + attrs + '>' +\n    '<w:r><w:t>Karl Sanchez</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> A0003707</w:t></w:r>' +\n    '<w:r><w:t>4</w:t></w:r>' +\n    '</w:p>' +\n    PKG_CLOSE;\n  karlPara.insertOoxml(xml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the commit's edits:\n# 1. \"WebDev Inc members:\" paragraph - drop the spell-check proofErr\n#    markers and merge the \"WebDev\" / \" Inc members:\" runs into a\n#    single run (keeping the bold formatting).\n# 2. \"Karl Sanchez\" paragraph - append a student ID, split across two\n#    extra runs: \" A0003707\" and \"4\" (reads as \"A00037074\").\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaOpenTagAttrs($range) {\n    $xml = $range.WordOpenXML\n    if ($xml -match '<w:p\\b([^>]*)>') {\n        return $matches[1]\n    }\n    return \"\"\n}\n\n$pkgOpen = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>'\n$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$membersPara = $null\n$karlPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if (($null -eq $membersPara) -and ($t -like \"*WebDev*\") -and ($t -like \"*Inc members*\")) {\n        $membersPara = $p\n    }\n    if (($null -eq $karlPara) -and ($t -like \"*Karl Sanchez*\")) {\n        $karlPara = $p\n    }\n}\n\nif ($null -ne $membersPara) {\n    $attrs = Get-ParaOpenTagAttrs $membersPara.Range\n    $xml = $pkgOpen + \"<w:p$attrs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>WebDev Inc members:</w:t></w:r></w:p>\" + $pkgClose\n    $membersPara.Range.InsertXML($xml)\n}\n\nif ($null -ne $karlPara) {\n    $attrs = Get-ParaOpenTagAttrs $karlPara.Range\n    $xml = $pkgOpen + \"<w:p$attrs><w:r><w:t>Karl Sanchez</w:t></w:r><w:r><w:t xml:space=`\"preserve`\"> A0003707</w:t></w:r><w:r><w:t>4</w:t></w:r></w:p>\" + $pkgClose\n    $karlPara.Range.InsertXML($xml)\n}\n"}
